$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.381.83'
$ws.Range("E2").Value = '  +9.11%  '
$ws.Range("D3").Value = '1.674.57'
$ws.Range("E3").Value = '  +4.63%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.000'
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '305.71'
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3683'
$ws.Range("E7").Value = '  +0.94%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3419'
$ws.Range("E8").Value = '  +0.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.48'
$ws.Range("E9").Value = '  +13.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.157'
$ws.Range("E10").Value = '  +2.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07215'
$ws.Range("E11").Value = '  +2.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.129'
$ws.Range("E13").Value = '  +4.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.02'
$ws.Range("E14").Value = '  +2.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.703'
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("D16").Value = '1.674.63'
$ws.Range("E16").Value = '  +4.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001097'
$ws.Range("E17").Value = '  +2.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9995'
$ws.Range("E18").Value = '  +0.54%  '
$ws.Range("E19").Value = '  +0.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '80.30'
$ws.Range("E20").Value = '  +3.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.40'
$ws.Range("E21").Value = '  +3.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.087'
$ws.Range("E22").Value = '  +1.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.25'
$ws.Range("E23").Value = '  +4.46%  '
$ws.Range("D24").Value = '24.349.10'
$ws.Range("E24").Value = '  +8.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.445'
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.644'
$ws.Range("E26").Value = '  +4.29%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.94'
$ws.Range("E27").Value = '  +1.95%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.39'
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D29").Value = '1.863.83'
$ws.Range("E29").Value = '  +5.37%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.78'
$ws.Range("E30").Value = '  +4.20%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.261'
$ws.Range("E31").Value = '  +3.58%  '
$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.050'
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9725'
$ws.Range("E33").Value = '  +4.29%  '
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08431'
$ws.Range("E34").Value = '  +2.92%  '
$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.674'
$ws.Range("E35").Value = '  -0.60%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.26'
$ws.Range("E36").Value = '  +3.73%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06382'
$ws.Range("E37").Value = '  +6.05%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.296'
$ws.Range("E38").Value = '  +2.44%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02316'
$ws.Range("E39").Value = '  +5.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.632'
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.238'
$ws.Range("E41").Value = '  -0.53%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2081'
$ws.Range("E42").Value = '  +3.86%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6065'
$ws.Range("E43").Value = '  +3.65%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  +0.60%  '
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.743'
$ws.Range("E45").Value = '  -1.80%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.96'
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5860'
$ws.Range("E47").Value = '  +3.98%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.55'
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.011'
$ws.Range("E49").Value = '  +3.00%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07137'
$ws.Range("E50").Value = '  +5.11%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '75.72'
$ws.Range("E51").Value = '  +3.63%  '
